$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking "crypto ticker" refresh: price (D) and 1h-volume (E) columns
# are plain text in the source feed (note multi-dot "26.758.85"-style prices),
# so any D-column value that Excel could mistake for a real number needs its
# cell pinned to Text format first to stop auto-coercion into a Double.

$ws.Range("D2").Value = "26.758.85"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.603.10"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.02"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "1.828.29"
$ws.Range("D13").Value = "1.601.37"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.02"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "0.0₃0739"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "209.77"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.14"
$ws.Range("E20").Value = "  +1.68%  "
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.30"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.24"
$ws.Range("E22").Value = "  -4.26%  "
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.60"
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.35"
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0509"
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("D33").Value = "1.289.98"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("E36").Value = "  -2.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.20"
$ws.Range("E37").Value = "  +13.12%  "
$ws.Range("D37").Style = "Normal"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.835"
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.43"
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("D40").Style = "Normal"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.781"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.03"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "1.740.17"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("D47").Value = "0.0₆0104"
$ws.Range("E47").Value = "  -2.53%  "
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0518"
$ws.Range("E49").Value = "  +1.78%  "
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.45"
$ws.Range("E51").Value = "  +1.02%  "
$ws.Range("D51").Style = "Normal"
